# Generate Report for handback
#
# The e716d8a8-7d6e-492c-b3ce-9abf02466350.md localization file has been
# handed back (in sync with en-US). Update the status on all three sheets
# and record the new handback timestamps for zh-cn / de-de.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = "Handed back: in sync with en-US"
$zhcn.Range("G3").Value = "2016-01-13 11:37:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = "Handed back: in sync with en-US"
$dede.Range("G3").Value = "2016-01-13 11:38:15"
